$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Add the new "local/lambda" / "local" values in row 6 (columns A/B)
# -----------------------------------------------------------------
$ws.Range("A6").Value = "local/lambda"
$ws.Range("B6").Value = "local"

# -----------------------------------------------------------------
# 2) Update the B21 link target/text from the old development1 URL
#    to https://www.advantageclub.co/ (same target as H9).
#    The engine only supports adding hyperlinks (not editing them in
#    place), and Worksheet.Hyperlinks.Delete() clears every hyperlink
#    on the sheet, so the reliable way to change just one hyperlink is
#    to drop them all and recreate them with the updated target for B21.
# -----------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:dheerajc@advantageclub.in", "", "", "dheerajc@advantageclub.in")
$ws.Hyperlinks.Add($ws.Range("B12"), "mailto:Dheeraj@4321", "", "", "Dheeraj@4321")
$ws.Hyperlinks.Add($ws.Range("H3"), "https://development1.advantageclub.co/in/rewards/home", "", "", "https://development1.advantageclub.co/in/rewards/home")
$ws.Hyperlinks.Add($ws.Range("H2"), "https://development1.advantageclub.co/", "", "https://development1.advantageclub.co/", "https://development1.advantageclub.co/")
$ws.Hyperlinks.Add($ws.Range("H5"), "https://celio.advantageclub.co/in/pages/rewards_home", "", "", "https://celio.advantageclub.co/in/pages/rewards_home")
$ws.Hyperlinks.Add($ws.Range("H4"), "https://celio.advantageclub.co/login", "", "", "https://celio.advantageclub.co/login")
$ws.Hyperlinks.Add($ws.Range("H6"), "https://secure.advantageclub.co/", "", "https://secure.advantageclub.co/", "https://secure.advantageclub.co/")
$ws.Hyperlinks.Add($ws.Range("H7"), "https://secure.workadvantage.in/rewards/home", "", "https://secure.workadvantage.in/rewards/home", "https://secure.workadvantage.in/rewards/home")
$ws.Hyperlinks.Add($ws.Range("H8"), "https://staging7.advantageclub.co/in/rewards/home", "", "https://staging7.advantageclub.co/in/rewards/home", "https://staging7.advantageclub.co/in/rewards/home")
$ws.Hyperlinks.Add($ws.Range("H9"), "https://www.advantageclub.co/", "", "https://www.advantageclub.co/", "https://www.advantageclub.co/")
$ws.Hyperlinks.Add($ws.Range("H11"), "https://staging7.advantageclub.co/in/login", "", "https://staging7.advantageclub.co/in/login", "https://staging7.advantageclub.co/in/login")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://development1.advantageclub.co/in/rewards/home", "", "", "https://development1.advantageclub.co/in/rewards/home")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://development1.advantageclub.co/", "", "https://development1.advantageclub.co/", "https://development1.advantageclub.co/")
$ws.Hyperlinks.Add($ws.Range("B21"), "https://www.advantageclub.co/", "", "https://www.advantageclub.co/", "https://www.advantageclub.co/")

# Re-adding hyperlinks resets the affected cells to a freshly minted
# "Hyperlink" style variant; restore the original named style so the
# cell formatting matches what it was before (style index unaffected).
$ws.Range("B11").Style = "Hyperlink"
$ws.Range("B12").Style = "Hyperlink"
$ws.Range("H3").Style = "Hyperlink"
$ws.Range("H2").Style = "Hyperlink"
$ws.Range("H5").Style = "Hyperlink"
$ws.Range("H4").Style = "Hyperlink"
$ws.Range("H6").Style = "Hyperlink"
$ws.Range("H7").Style = "Hyperlink"
$ws.Range("H8").Style = "Hyperlink"
$ws.Range("H9").Style = "Hyperlink"
$ws.Range("H11").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B21").Style = "Hyperlink"

# -----------------------------------------------------------------
# 3) Move the active selection to C17, matching the saved view state.
# -----------------------------------------------------------------
$ws.Range("C17").Select()
